$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add damage reduction stats (esprit, endurance, robustesse) for level 1
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 3

# Increase enemies damage (intelligence) for level 1
$ws.Range("F2").Value = 2

# Update the active selection to J3
$ws.Range("J3").Select()
